$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style of a "normal" (unstyled) data cell, used to strip the extra
# number-format style that gets attached when forcing text-typed numeric strings.
$normalStyle = $ws.Range("B2").Style

$ws.Range("D2").Value2 = '27.988.57'
$ws.Range("E2").Value2 = '  -5.16%  '

$ws.Range("D3").Value2 = '1.824.91'
$ws.Range("E3").Value2 = '  -4.00%  '

$ws.Range("E4").Value2 = '  -0.26%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value2 = '328.84'
$cell.Style = $normalStyle
$ws.Range("E5").Value2 = '  -2.81%  '

$ws.Range("E6").Value2 = '  -0.29%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value2 = '0.4631'
$cell.Style = $normalStyle
$ws.Range("E7").Value2 = '  -2.55%  '

$ws.Range("E8").Value2 = '  -3.48%  '

$ws.Range("B9").Value2 = 'Dogecoin'
$ws.Range("C9").Value2 = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value2 = '0.07871'
$cell.Style = $normalStyle
$ws.Range("E9").Value2 = '  -2.15%  '

$ws.Range("B10").Value2 = 'Polygon'
$ws.Range("C10").Value2 = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value2 = '0.9590'
$cell.Style = $normalStyle
$ws.Range("E10").Value2 = '  -3.23%  '

$ws.Range("B11").Value2 = 'Solana'
$ws.Range("C11").Value2 = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value2 = '21.87'
$cell.Style = $normalStyle
$ws.Range("E11").Value2 = '  -5.75%  '

$ws.Range("B12").Value2 = 'WrappedEther'
$ws.Range("C12").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value2 = '1.802.17'
$ws.Range("E12").Value2 = '  -6.12%  '

$ws.Range("B13").Value2 = 'Polkadot'
$ws.Range("C13").Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value2 = '5.651'
$cell.Style = $normalStyle
$ws.Range("E13").Value2 = '  -4.79%  '

$ws.Range("B14").Value2 = 'Chainlink'
$ws.Range("C14").Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value2 = '6.884'
$cell.Style = $normalStyle
$ws.Range("E14").Value2 = '  -3.05%  '

$ws.Range("B15").Value2 = 'TRON'
$ws.Range("C15").Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value2 = '0.06833'
$cell.Style = $normalStyle
$ws.Range("E15").Value2 = '  +0.25%  '

$ws.Range("B16").Value2 = 'BinanceUSD'
$ws.Range("C16").Value2 = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value2 = '1.001'
$cell.Style = $normalStyle
$ws.Range("E16").Value2 = '  -0.35%  '

$ws.Range("B17").Value2 = 'Litecoin'
$ws.Range("C17").Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value2 = '86.41'
$cell.Style = $normalStyle
$ws.Range("E17").Value2 = '  -3.01%  '

$ws.Range("B18").Value2 = 'ShibaInu'
$ws.Range("C18").Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value2 = '0.000009965'
$cell.Style = $normalStyle
$ws.Range("E18").Value2 = '  -2.31%  '

$ws.Range("B19").Value2 = 'Avalanche'
$ws.Range("C19").Value2 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value2 = '16.64'
$cell.Style = $normalStyle
$ws.Range("E19").Value2 = '  -4.01%  '

$ws.Range("B20").Value2 = 'Dai'
$ws.Range("C20").Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value2 = '1.0000'
$cell.Style = $normalStyle
$ws.Range("E20").Value2 = '  -0.32%  '

$ws.Range("B21").Value2 = 'WrappedBTC'
$ws.Range("C21").Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value2 = '28.021.79'
$ws.Range("E21").Value2 = '  -5.11%  '

$ws.Range("B22").Value2 = 'Uniswap'
$ws.Range("C22").Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value2 = '5.318'
$cell.Style = $normalStyle
$ws.Range("E22").Value2 = '  -3.48%  '

$ws.Range("B23").Value2 = 'Cosmos'
$ws.Range("C23").Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value2 = '10.98'
$cell.Style = $normalStyle
$ws.Range("E23").Value2 = '  -5.44%  '

$ws.Range("B24").Value2 = 'Toncoin'
$ws.Range("C24").Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value2 = '2.091'
$cell.Style = $normalStyle
$ws.Range("E24").Value2 = '  -2.78%  '

$ws.Range("B25").Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value2 = '2.056.70'
$ws.Range("E25").Value2 = '  -4.10%  '

$ws.Range("B26").Value2 = 'Monero'
$ws.Range("C26").Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value2 = '152.12'
$cell.Style = $normalStyle
$ws.Range("E26").Value2 = '  -3.19%  '

$ws.Range("B27").Value2 = 'EthereumClassic'
$ws.Range("C27").Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value2 = '19.14'
$cell.Style = $normalStyle
$ws.Range("E27").Value2 = '  -2.53%  '

$ws.Range("B28").Value2 = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value2 = '5.740'
$cell.Style = $normalStyle
$ws.Range("E28").Value2 = '  -11.43%  '

$ws.Range("B29").Value2 = 'LidoDAOToken'
$ws.Range("C29").Value2 = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value2 = '1.966'
$cell.Style = $normalStyle
$ws.Range("E29").Value2 = '  -4.43%  '

$ws.Range("B30").Value2 = 'BitcoinCash'
$ws.Range("C30").Value2 = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value2 = '116.62'
$cell.Style = $normalStyle
$ws.Range("E30").Value2 = '  -2.09%  '

$ws.Range("B31").Value2 = 'ImmutableX'
$ws.Range("C31").Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value2 = '0.9368'
$cell.Style = $normalStyle
$ws.Range("E31").Value2 = '  -5.84%  '

$ws.Range("B32").Value2 = 'Stellar'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value2 = '0.09215'
$cell.Style = $normalStyle
$ws.Range("E32").Value2 = '  -3.43%  '

$ws.Range("B33").Value2 = 'Filecoin'
$ws.Range("C33").Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value2 = '5.283'
$cell.Style = $normalStyle
$ws.Range("E33").Value2 = '  -3.52%  '

$ws.Range("B34").Value2 = 'ARBITRUM'
$ws.Range("C34").Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value2 = '1.316'
$cell.Style = $normalStyle
$ws.Range("E34").Value2 = '  -5.32%  '

$ws.Range("B35").Value2 = 'HuobiToken'
$ws.Range("C35").Value2 = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value2 = '3.344'
$cell.Style = $normalStyle
$ws.Range("E35").Value2 = '  -5.21%  '

$ws.Range("B36").Value2 = 'Hedera'
$ws.Range("C36").Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value2 = '0.05928'
$cell.Style = $normalStyle
$ws.Range("E36").Value2 = '  -7.18%  '

$ws.Range("B37").Value2 = 'VeChain'
$ws.Range("C37").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value2 = '0.02139'
$cell.Style = $normalStyle
$ws.Range("E37").Value2 = '  -4.79%  '

$ws.Range("B38").Value2 = 'TrustWalletToken'
$ws.Range("C38").Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value2 = '1.146'
$cell.Style = $normalStyle
$ws.Range("E38").Value2 = '  -4.21%  '

$ws.Range("B39").Value2 = 'Frax'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value2 = '1.0000'
$cell.Style = $normalStyle
$ws.Range("E39").Value2 = '  -0.35%  '

$ws.Range("B40").Value2 = 'FraxShare'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value2 = '7.603'
$cell.Style = $normalStyle
$ws.Range("E40").Value2 = '  -1.61%  '

$ws.Range("B41").Value2 = 'TheSandbox'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value2 = '0.5570'
$cell.Style = $normalStyle
$ws.Range("E41").Value2 = '  -4.30%  '

$ws.Range("B42").Value2 = 'Aptos'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value2 = '9.895'
$cell.Style = $normalStyle
$ws.Range("E42").Value2 = '  -6.33%  '

$ws.Range("B43").Value2 = 'Algorand'
$ws.Range("C43").Value2 = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value2 = '0.1763'
$cell.Style = $normalStyle
$ws.Range("E43").Value2 = '  -3.05%  '

$ws.Range("B44").Value2 = 'WEMIXToken'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value2 = '1.223'
$cell.Style = $normalStyle
$ws.Range("E44").Value2 = '  -3.35%  '

$ws.Range("B45").Value2 = 'RenderToken'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value2 = '2.218'
$cell.Style = $normalStyle
$ws.Range("E45").Value2 = '  -8.58%  '

$ws.Range("B46").Value2 = 'EnergySwap'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value2 = '11.58'
$cell.Style = $normalStyle
$ws.Range("E46").Value2 = '  -4.29%  '

$ws.Range("B47").Value2 = 'Decentraland'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value2 = '0.5252'
$cell.Style = $normalStyle
$ws.Range("E47").Value2 = '  -4.35%  '

$ws.Range("B48").Value2 = 'Cronos'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value2 = '0.07002'
$cell.Style = $normalStyle
$ws.Range("E48").Value2 = '  -4.91%  '

$ws.Range("B49").Value2 = 'NEARProtocol'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value2 = '1.823'
$cell.Style = $normalStyle
$ws.Range("E49").Value2 = '  -6.79%  '

$ws.Range("B50").Value2 = 'Quant'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value2 = '111.17'
$cell.Style = $normalStyle
$ws.Range("E50").Value2 = '  -4.49%  '

$ws.Range("B51").Value2 = 'PaxDollar'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value2 = '0.9998'
$cell.Style = $normalStyle
$ws.Range("E51").Value2 = '  -0.34%  '
